# ESS11 edition 4.0: new country row data added to tab_weight_sums.
# A new row for "EE" (Estonia) is inserted before the existing row 114,
# shifting the old rows 114-135 down to 115-136, and a new row for "UA"
# (Ukraine) is appended after the old last row (now row 136), becoming
# row 137.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 114 (currently holding "ES"),
# pushing rows 114-135 down to 115-136.
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with the "EE" (Estonia) data.
$ws.Cells.Item(114, 1).Value = "R11"
$ws.Cells.Item(114, 2).Value = "EE"
$ws.Cells.Item(114, 3).Value = 1293
$ws.Cells.Item(114, 4).Value = 1292.99997454882
$ws.Cells.Item(114, 5).Value = 1292.99996763468
$ws.Cells.Item(114, 6).Value = 1141964.9771601
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 0

# Append the new "UA" (Ukraine) row after the old last row (shifted to 136).
$ws.Cells.Item(137, 1).Value = "R11"
$ws.Cells.Item(137, 2).Value = "UA"
$ws.Cells.Item(137, 3).Value = 2661
$ws.Cells.Item(137, 4).Value = 2661.00006902218
$ws.Cells.Item(137, 5).Value = 2661.00006902218
$ws.Cells.Item(137, 6).Value = 34877812.7890825
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 0
